# Update todo list and player attribute script

$wb = $excel.ActiveWorkbook

# --- 1. Append new Todo items to the "Todo " sheet ---
$todo = $wb.Worksheets.Item("Todo ")

$todo.Range("A34").Value = "Adjust camera"
$todo.Range("D34").Value = "Fish"
$todo.Range("F34").Value = "18 Jan"

$todo.Range("A35").Value = "Migreate light"
$todo.Range("C35").Value = "Forest"
$todo.Range("D35").Value = "Fish"
$todo.Range("F35").Value = "18 Jan"

$todo.Range("A36").Value = "Log console"
$todo.Range("D36").Value = "Fish"
$todo.Range("F36").Value = "18 Jan"

# --- 2. Add a new "Remark" worksheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$remark = $wb.Worksheets.Add($null, $lastSheet)
$remark.Name = "Remark"

$remark.Range("A1").Value = "Item"
$remark.Range("B1").Value = "0 - 99"

$remark.Range("A2").Value = "Weapon"
$remark.Range("B2").Value = "100 - 199"

$remark.Range("A3").Value = "Armor"
$remark.Range("B3").Value = "200 - 299"

# --- 3. Fix up the selections that Excel leaves behind ---
# Move the Todo sheet's selection cursor past the newly added rows...
[void]$todo.Range("A37").Select()

# ...then make Remark the active/selected tab again, with B3 selected
[void]$remark.Activate()
[void]$remark.Range("B3").Select()
